# Update countries & provincias Spain
# - Refresh the "last updated" timestamp
# - Update India's and Mongolia's case counts
# - Re-sort a handful of small territories (Sudan del Sur, Bonaire, Islas
#   Virgenes Britanicas, Santo Tome y Principe, San Pedro y Miquelon, Yemen)
#   and refresh their case counts

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 05:52"

# --- India (row 19): Casos totales, Nuevos casos, Recuperados ---
$ws.Cells.Item(19, 2).Value = 23077
$ws.Cells.Item(19, 3).Value = 38
$ws.Cells.Item(19, 5).Value = 17344

# --- Mongolia (row 172): Casos totales, Nuevos casos, Casos activos ---
$ws.Cells.Item(172, 2).Value = 36
$ws.Cells.Item(172, 3).Value = 1
$ws.Cells.Item(172, 4).Value = 9

# --- Small territories (rows 210-216) reshuffled with refreshed counts ---
# Row 210: Sudan del Sur
$ws.Cells.Item(210, 1).Value = "Sudan del Sur"
$ws.Cells.Item(210, 2).Value = 5
$ws.Cells.Item(210, 3).Value = 1
$ws.Cells.Item(210, 4).Value = 0
$ws.Cells.Item(210, 5).Value = 5
$ws.Cells.Item(210, 6).Value = 0
$ws.Cells.Item(210, 7).Value = 0
$ws.Cells.Item(210, 8).Value = 0

# Row 211: Bonaire, San Eustaquio y Saba
$ws.Cells.Item(211, 1).Value = "Bonaire, San Eustaquio y Saba"
$ws.Cells.Item(211, 2).Value = 5
$ws.Cells.Item(211, 3).Value = 0
$ws.Cells.Item(211, 4).Value = 0
$ws.Cells.Item(211, 5).Value = 5
$ws.Cells.Item(211, 6).Value = 0
$ws.Cells.Item(211, 7).Value = 0
$ws.Cells.Item(211, 8).Value = 0

# Row 212: Islas Virgenes Britanicas
$ws.Cells.Item(212, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(212, 2).Value = 5
$ws.Cells.Item(212, 3).Value = 0
$ws.Cells.Item(212, 4).Value = 3
$ws.Cells.Item(212, 5).Value = 1
$ws.Cells.Item(212, 6).Value = 0
$ws.Cells.Item(212, 7).Value = 0
$ws.Cells.Item(212, 8).Value = 1

# Row 213: Santo Tome y Principe
$ws.Cells.Item(213, 1).Value = "Santo Tome y Principe"
$ws.Cells.Item(213, 2).Value = 4
$ws.Cells.Item(213, 3).Value = 0
$ws.Cells.Item(213, 4).Value = 0
$ws.Cells.Item(213, 5).Value = 4
$ws.Cells.Item(213, 6).Value = 0
$ws.Cells.Item(213, 7).Value = 0
$ws.Cells.Item(213, 8).Value = 0

# Row 214: Anguila (unchanged values, kept for completeness)
$ws.Cells.Item(214, 1).Value = "Anguila"
$ws.Cells.Item(214, 2).Value = 3
$ws.Cells.Item(214, 3).Value = 0
$ws.Cells.Item(214, 4).Value = 1
$ws.Cells.Item(214, 5).Value = 2
$ws.Cells.Item(214, 6).Value = 0
$ws.Cells.Item(214, 7).Value = 0
$ws.Cells.Item(214, 8).Value = 0

# Row 215: San Pedro y Miquelon
$ws.Cells.Item(215, 1).Value = "San Pedro y Miquelon"
$ws.Cells.Item(215, 2).Value = 1
$ws.Cells.Item(215, 3).Value = 0
$ws.Cells.Item(215, 4).Value = 0
$ws.Cells.Item(215, 5).Value = 1
$ws.Cells.Item(215, 6).Value = 0
$ws.Cells.Item(215, 7).Value = 0
$ws.Cells.Item(215, 8).Value = 0

# Row 216: Yemen
$ws.Cells.Item(216, 1).Value = "Yemen"
$ws.Cells.Item(216, 2).Value = 1
$ws.Cells.Item(216, 3).Value = 0
$ws.Cells.Item(216, 4).Value = 1
$ws.Cells.Item(216, 5).Value = 0
$ws.Cells.Item(216, 6).Value = 0
$ws.Cells.Item(216, 7).Value = 0
$ws.Cells.Item(216, 8).Value = 0
